# no-op test
$wb = $excel.ActiveWorkbook
Write-Host "Sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
